$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The author removed the now-unnecessary footnote/legend rows (20-26) that
# sat below the main Schulferien table (commit: "unnötige Zeilen gelöscht").
# Clear the cell contents (values) while leaving row heights / cell
# formatting untouched, mirroring an Excel "Clear Contents" on that block.
$ws.Range("A20:G26").ClearContents()

# Match the resulting selection left behind in the sheet (A18, with the
# block A18:H30 highlighted) after the rows were cleared.
$ws.Range("A18:H30").Select()
